$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Change 1: split the run
#   "bon de la rougir au foeu pour brusler les immundices"
# into four runs:
#   "bon de la rougir au foeu pour brusler les "  (plain)
#   "<m>"                                          (Courier New, blue, 9pt)
#   "immundices"                                   (plain)
#   "</m>"                                         (Courier New, blue, 9pt)
# ----------------------------------------------------------------------

$rng1 = $d.Content
$found1 = $rng1.Find.Execute("immundices", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $startPos = $rng1.Start
    $endPos = $rng1.End

    # Insert "</m>" right after "immundices" and format it.
    $afterPoint = $d.Range($endPos, $endPos)
    $afterPoint.InsertAfter("</m>")
    $closeTag = $d.Range($endPos, $endPos + 4)
    $closeTag.Font.Name = "Courier New"
    $closeTag.Font.Color = 16711680
    $closeTag.Font.Size = 9

    # Insert "<m>" right before "immundices" and format it.
    $beforePoint = $d.Range($startPos, $startPos)
    $beforePoint.InsertBefore("<m>")
    $openTag = $d.Range($startPos, $startPos + 3)
    $openTag.Font.Name = "Courier New"
    $openTag.Font.Color = 16711680
    $openTag.Font.Size = 9
}

# ----------------------------------------------------------------------
# Change 2: in the paragraph "bonne <m>limaille</m> ira au fons."
#   - delete the leading run "bonne "
#   - change the "limaille" run's text to "bonne limaille"
# ----------------------------------------------------------------------

$rng2 = $d.Content
$found2 = $rng2.Find.Execute("bonne <m>limaille", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    # "bonne " is the first 6 characters of this match.
    $bonneRange = $d.Range($rng2.Start, $rng2.Start + 6)
    $bonneRange.Delete()
}

$rng3 = $d.Content
$found3 = $rng3.Find.Execute("limaille</m> ira au fons", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found3) {
    # "limaille" is the first 8 characters of this match.
    $limRange = $d.Range($rng3.Start, $rng3.Start + 8)
    $limRange.Text = "bonne limaille"
}
